$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Contact Person's Email (row 10): N/A -> Daniel_CHEONG@jtc.gov.sg
$ws.Range("B10").Value = "Daniel_CHEONG@jtc.gov.sg"

# Contact Person's Tel Num (row 11): N/A -> 0 (numeric)
$ws.Range("B11").Value = 0

# Contact Person's Fax Num (row 12): N/A -> blank (cell removed)
$ws.Range("B12").ClearContents()

# Contact Person's Address (row 13): N/A -> The JTC Summit, 8 Jurong Town Hall Road, Singapore 609434
$ws.Range("B13").Value = "The JTC Summit, 8 Jurong Town Hall Road, Singapore 609434"
